# Gestion d'accès au batiment.pptx -- "Mise à jour planches"
#
# 1. Sommaire (slide 2): add a 5th bullet "Répartition des tâches".
# 2. Composants/Matériels (slide 6):
#      - "...présente dans la banque." -> "...présente dans le bâtiment."
#      - "Lecteur RFID" -> "M5Stack"
# 3. Two new slides appended at the end (Title+Content layout):
#      slide 7: "5. Répartition des tâches"
#      slide 8: "6. Tâches réalisés"

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Slide 2 ("Sommaire") -- append new agenda bullet
# ---------------------------------------------------------------------
$sommaire = $p.Slides.Item(2)
$sommaireBody = $sommaire.Shapes.Item(2).TextFrame.TextRange
$lastPara = $sommaireBody.Paragraphs($sommaireBody.Paragraphs().Count, 1)
$newBullet = $lastPara.InsertAfter("`rRépartition des tâches")
$newBullet.LanguageID = "fr-FR"

# ---------------------------------------------------------------------
# 2. Slide 6 ("Composants/Matériels") -- text tweaks
# ---------------------------------------------------------------------
$composants = $p.Slides.Item(6)
$composantsBody = $composants.Shapes.Item(2).TextFrame.TextRange

# "...présente dans la banque." -> "...présente dans le bâtiment."
$tvPara = $composantsBody.Paragraphs(1, 1)
$tailStart = $tvPara.Text.IndexOf("la banque.") + 1
$tail = $tvPara.Characters($tailStart, 10)
$tail.Text = "le bâtiment."
$tail.LanguageID = "fr-FR"

# "Lecteur RFID" -> "M5Stack"
$lecteurPara = $composantsBody.Paragraphs(6, 1)
$lecteurPara.Text = "M5Stack"
$lecteurPara.LanguageID = "fr-FR"

# ---------------------------------------------------------------------
# 3. New slide 7 -- "5. Répartition des tâches"
# ---------------------------------------------------------------------
$slide7 = $p.Slides.Add($p.Slides.Count + 1, 2)

$s7Title = $slide7.Shapes.Item(1).TextFrame.TextRange
$s7Title.Text = "5. Répartition des tâches"
$s7Title.LanguageID = "fr-FR"

$s7Body = $slide7.Shapes.Item(2).TextFrame.TextRange
$s7Body.Text = "Saxemard:"
$s7Body.LanguageID = "fr-FR"

$r = $s7Body.InsertAfter(" ")
$r.LanguageID = "fr-FR"
$r = $s7Body.InsertAfter("Mise en place d’un système de création de créneau sur un agenda Google, service qui va avec, mise en place de la connexion au broker MQTT pour l’envoi des données (Code, id…)")
$r.LanguageID = "fr-FR"

$r = $s7Body.InsertAfter("`rPillar")
$r.LanguageID = "fr-FR"
$r = $s7Body.InsertAfter(" : ")
$r.LanguageID = "fr-FR"
$r = $s7Body.InsertAfter("NodeRed")
$r.LanguageID = "fr-FR"
$r = $s7Body.InsertAfter(" ")
$r.LanguageID = "fr-FR"

$r = $s7Body.InsertAfter("`rPascucci")
$r.LanguageID = "fr-FR"
$r = $s7Body.InsertAfter(" : Mise en place de la base de données pour recueillir les informations provenant du M5Stack")
$r.LanguageID = "fr-FR"

# ---------------------------------------------------------------------
# 4. New slide 8 -- "6. Tâches réalisés"
# ---------------------------------------------------------------------
$slide8 = $p.Slides.Add($p.Slides.Count + 1, 2)

$s8Title = $slide8.Shapes.Item(1).TextFrame.TextRange
$s8Title.Text = "6. Tâches réalisés"
$s8Title.LanguageID = "fr-FR"

$s8Body = $slide8.Shapes.Item(2).TextFrame.TextRange
$s8Body.Text = "Saxemard : Dessin de l’architecture et répartition des tâches"
$s8Body.LanguageID = "fr-FR"
